$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data (Date serial, Method, ElapsedMs, wordCount, sentenceCount,
# posWordCount, negWordCount, posWordPercentage, negWordPercentage,
# positivePhraseCount, negativePhraseCount, posPhrasePercentage, negPhrasePercentage)
$data = @(
    @(42602.576018518521, "Noun", 3157, 326, 31, 7, 15, 31, 68, 0, 3, 0, 99),
    @(42602.576863425929, "Noun", 3028, 326, 31, 7, 15, 31, 68, 0, 3, 0, 99),
    @(42602.577187499999, "Noun", 2192, 326, 31, 4, 17, 19, 80, 0, 4, 0, 100)
)

$startRow = 27
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]

    # Copy the date-formatted style from an existing row (row 2) onto column A
    # of the new row so the same cell style index gets reused.
    $ws.Range("A2").Copy()
    $ws.Range("A" + $row).PasteSpecial(-4122)

    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $rowData[$c]
    }
}
